# Add data for 2022-01-01: update the sheet name, the December label, and
# the December / Total row values for the 2021 carjacking-by-month-yoy sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2021-12-24"

# Update the December row label.
$ws.Range("A13").Value = "December (through 12-24)"

# Update December row (row 13) values for years 2015-2021 (columns B-H).
$ws.Range("B13").Value = 34
$ws.Range("C13").Value = 77
$ws.Range("D13").Value = 94
$ws.Range("E13").Value = 56
$ws.Range("F13").Value = 50
$ws.Range("G13").Value = 116
$ws.Range("H13").Value = 156

# Update Total row (row 14) values for years 2015-2021 (columns B-H).
$ws.Range("B14").Value = 325
$ws.Range("C14").Value = 640
$ws.Range("D14").Value = 915
$ws.Range("E14").Value = 738
$ws.Range("F14").Value = 584
$ws.Range("G14").Value = 1380
$ws.Range("H14").Value = 1799
